$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input_control")

$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
